# Applies the cryptos-list refresh described by the commit:
# "Updated cryptos list on Sun Sep  8 21:22:48 UTC 2024 with GitHub Actions"
#
# Column D (Price) and E (Volume(1h)) hold values that LOOK numeric
# (e.g. "498.99", "0.996", "  +0.84%  ") but must stay plain TEXT, exactly
# like the original inline strings in the sheet. Excel's COM layer will
# silently coerce a bare numeric-looking string into a real Number when we
# assign .Value, so for any cell whose new text would parse as a number we
# first force the cell to the Text number format ("@") - this is the COM
# equivalent of pre-formatting the cell as Text before typing the value, and
# keeps values such as "245.40" or "0.0000130" from losing their exact
# formatting or turning into 245.4 / 1.3E-05.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.418.08'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '2.273.76'
$ws.Range("E3").Value = '  +1.26%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '498.99'
$ws.Range("E5").Value = '  +1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.73'
$ws.Range("E6").Value = '  +1.77%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.527'
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0959'
$ws.Range("E9").Value = '  +0.97%  '
$ws.Range("E10").Value = '  +0.95%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.337'
$ws.Range("E11").Value = '  +3.61%  '
$ws.Range("E12").Value = '  +5.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.13'
$ws.Range("E13").Value = '  +6.52%  '
$ws.Range("D14").Value = '2.674.70'
$ws.Range("E14").Value = '  +1.19%  '
$ws.Range("D15").Value = '54.397.17'
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000130'
$ws.Range("E16").Value = '  +1.18%  '
$ws.Range("D17").Value = '2.279.66'
$ws.Range("E17").Value = '  +0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.30'
$ws.Range("E18").Value = '  +2.87%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.16'
$ws.Range("E19").Value = '  +1.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '305.09'
$ws.Range("E20").Value = '  +1.93%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.35'
$ws.Range("E21").Value = '  -0.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '60.70'
$ws.Range("E23").Value = '  -1.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.996'
$ws.Range("E24").Value = '  -2.08%  '
$ws.Range("E25").Value = '  +1.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.37'
$ws.Range("E26").Value = '  +4.79%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '174.32'
$ws.Range("E27").Value = '  +4.69%  '
$ws.Range("D28").Value = '0.0₃0711'
$ws.Range("E28").Value = '  +4.12%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.02'
$ws.Range("E29").Value = '  +3.08%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.61'
$ws.Range("E30").Value = '  +0.87%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.09'
$ws.Range("E31").Value = '  +2.35%  '
$ws.Range("E32").Value = '  +0.03%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.88'
$ws.Range("E33").Value = '  +1.41%  '
$ws.Range("E34").Value = '  +0.19%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.952'
$ws.Range("E35").Value = '  +5.47%  '
$ws.Range("E36").Value = '  +2.72%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.74'
$ws.Range("E37").Value = '  +1.74%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.377'
$ws.Range("E38").Value = '  +1.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.41'
$ws.Range("E39").Value = '  +0.90%  '
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.11'
$ws.Range("E40").Value = '  +3.67%  '
$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.39'
$ws.Range("E41").Value = '  +1.36%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '125.58'
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0494'
$ws.Range("E43").Value = '  +2.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0898'
$ws.Range("E44").Value = '  +1.50%  '
$ws.Range("B45").Value = 'Bittensor'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '245.40'
$ws.Range("E45").Value = '  +3.94%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.549'
$ws.Range("E46").Value = '  +1.60%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.376'
$ws.Range("E47").Value = '  +1.78%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0206'
$ws.Range("E48").Value = '  +2.03%  '
$ws.Range("E49").Value = '  +0.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '16.31'
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.53'
$ws.Range("E51").Value = '  +3.13%  '
